$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.079.82'
$ws.Range('D2').Style = 'Normal'

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.494.87'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value2 = '  -1.21%  '

$ws.Range('E4').Value2 = '  +0.06%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '570.06'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value2 = '  -1.09%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '164.78'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value2 = '  -2.87%  '

$ws.Range('E7').Value2 = '  -0.07%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.510'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value2 = '  -0.15%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.491.63'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value2 = '  -1.35%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.157'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value2 = '  -1.85%  '

$ws.Range('E11').Value2 = '  -0.67%  '

$ws.Range('E12').Value2 = '  +2.75%  '

$ws.Range('E13').Value2 = '  +0.89%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.938.03'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value2 = '  -1.63%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '68.896.03'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value2 = '  -1.89%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000174'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value2 = '  -2.86%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '24.58'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value2 = '  -1.91%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.499.78'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value2 = '  -0.88%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.22'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value2 = '  -2.63%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.55'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value2 = '  -0.88%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '346.73'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value2 = '  -2.50%  '

$ws.Range('E22').Value2 = '  -1.84%  '

$ws.Range('E23').Value2 = '  -0.38%  '

$ws.Range('E24').Value2 = '  +0.00%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '69.91'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value2 = '  +0.95%  '

$ws.Range('E26').Value2 = '  -4.69%  '

$ws.Range('E27').Value2 = '  -4.30%  '

$ws.Range('E29').Value2 = '  +0.45%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0879'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value2 = '  -3.62%  '

$ws.Range('E31').Value2 = '  -1.62%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '457.34'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value2 = '  -5.51%  '

$ws.Range('E33').Value2 = '  -7.58%  '

$ws.Range('E34').Value2 = '  -2.76%  '

$ws.Range('E35').Value2 = '  +0.03%  '

$ws.Range('B36').Value2 = 'Monero'
$ws.Range('C36').Value2 = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '157.19'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value2 = '  +0.81%  '

$ws.Range('B37').Value2 = 'Kaspa'
$ws.Range('C37').Value2 = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.115'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value2 = '  -1.05%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.00'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value2 = '  +0.59%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.34'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value2 = '  -1.52%  '

$ws.Range('E40').Value2 = '  -0.05%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.315'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value2 = '  -1.80%  '

$ws.Range('E42').Value2 = '  -2.04%  '

$ws.Range('E43').Value2 = '  -3.60%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '38.02'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value2 = '  -0.74%  '

$ws.Range('E45').Value2 = '  -8.85%  '

$ws.Range('E46').Value2 = '  -8.39%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '141.09'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value2 = '  -1.17%  '

$ws.Range('E48').Value2 = '  -2.53%  '

$ws.Range('E49').Value2 = '  -2.45%  '

$ws.Range('E50').Value2 = '  -0.20%  '

$ws.Range('E51').Value2 = '  -4.37%  '

